$wb = $excel.ActiveWorkbook

# Update "展览" sheet (sheet1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 116
$ws1.Range("F4").Value = 65

# Update "全部类型" sheet (sheet4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 116
$ws4.Range("F4").Value = 65
